$wb = $excel.ActiveWorkbook

# Update data values on "inputdisp" sheet
$wsInput = $wb.Worksheets.Item("inputdisp")
$wsInput.Range("B2").Value = 1.5
$wsInput.Range("B3").Value = 3

# Make "inputdisp" the active sheet (was "endofpipe"), and move the
# active-cell selection on that sheet to B4
$wsInput.Activate()
$wsInput.Range("B4").Select()
